$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old summary block (blank template row 23, sum[min] row 24, sum[h] row
# 25, sum[working weeks] row 26) all shifts down by one row to make room for
# a new data row at 23. Copy bottom-up, cell by cell (format + content), so
# a row is never clobbered before it has been copied from.
26..23 | ForEach-Object {
    $src = $_
    $dst = $_ + 1

    # Column D only matters for the blank template row (23 -> 24); rows
    # 24-26 never had a D cell to begin with.
    if ($src -eq 23) {
        $dCell = $ws.Range("D$dst")
        $dCell.Clear()
        $dCell.NumberFormat = $ws.Range("D$src").NumberFormat
    }

    # Column E: the blank template row (23 -> 24) needs hh:mm formatting;
    # the "sum [min]" / "sum [h]" / "sum [working weeks]" label rows just
    # need right-alignment (row 27 is brand new so needs it applied
    # explicitly; rows 25/26 already have it from their old content, but
    # setting it again is harmless and keeps things uniform).
    $eSrc = $ws.Range("E$src")
    $eDst = $ws.Range("E$dst")
    $eDst.Clear()
    if ($src -eq 23) {
        $eDst.NumberFormat = $eSrc.NumberFormat
    } else {
        $eDst.HorizontalAlignment = -4152
    }
    if ($eSrc.HasFormula) { $eDst.Formula = $eSrc.Formula } else { $eDst.Value = $eSrc.Value2 }

    # Column F: formatting differs row to row (integer for the blank
    # template/sum[min] rows, 2-decimal for the sum[h]/sum[working weeks]
    # rows), so always carry the source's number format down with it.
    $fSrc = $ws.Range("F$src")
    $fDst = $ws.Range("F$dst")
    $fDst.NumberFormat = $fSrc.NumberFormat
    if ($fSrc.HasFormula) { $fDst.Formula = $fSrc.Formula } else { $fDst.Value = $fSrc.Value2 }
}

# --- Row 23: new data row (overwrite what got copied down from the old
#     blank template row) ---
$ws.Range("A23").Value = 2014
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = 24
$ws.Range("D23").Value = 0.75
$ws.Range("E23").Value = 0.83333333333333337

# Extend the "time spent [min]" / "time spent [h]" formulas down through
# the new row 23, same as dragging their fill handle down one more row.
$ws.Range("F23").Formula = "=(E23-D23)*24*60"
$ws.Range("G23").Formula = "=F23/60"

# --- Fix up the formulas in the shifted-down summary rows so they point at
#     their new (shifted) source rows, since the literal copy above kept the
#     old row references. ---
$ws.Range("F25").Formula = "=SUM(F2:F24)"
$ws.Range("F26").Formula = "=F25/60"
$ws.Range("F27").Formula = "=F26/38.5"

$ws.Range("F23").Select()
